# "Ajuste de pum en carga de precio base y precio por rol"
#
# The PUM column (column C: "PUM" header + "Gramo a $X,6 pesos" labels)
# is removed entirely from the import template. The former column D
# ("Mostrar Descuento (1 Si, 0 No)" header + its 1/0 values and the
# style-only marker cells) shifts left into column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole PUM column (C). This removes the "PUM" header and the
# "Gramo a $X,6 pesos" text values, and shifts the old column D
# (Mostrar Descuento) left to become the new column C.
$ws.Columns.Item(3).Delete()

# Move the selection cursor to where it ended up after the edit.
$ws.Range("C12").Select()
